# Auto-generated Excel COM-interop script to apply scheduled market-data refresh
# Updates currentAveragePrice* / LevePrice* / LeveProfit* columns (H-N) across all 8 sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(28, 8).Value = 349.875
$ws.Cells.Item(28, 9).Value = 349.875
$ws.Cells.Item(28, 11).Value = 349.875
$ws.Cells.Item(28, 13).Value = 135.125
$ws.Cells.Item(58, 8).Value = 3141
$ws.Cells.Item(58, 9).Value = 176.25
$ws.Cells.Item(58, 10).Value = 15000
$ws.Cells.Item(58, 11).Value = 528.75
$ws.Cells.Item(58, 12).Value = 45000
$ws.Cells.Item(58, 13).Value = -378.75
$ws.Cells.Item(58, 14).Value = -45300
$ws.Cells.Item(80, 8).Value = 3597.6667
$ws.Cells.Item(80, 9).Value = 0
$ws.Cells.Item(80, 10).Value = 3597.6667
$ws.Cells.Item(80, 11).Value = 0
$ws.Cells.Item(80, 12).Value = 10793.0001
$ws.Cells.Item(80, 14).Value = -12789.0001
$ws.Cells.Item(83, 8).Value = 3597.6667
$ws.Cells.Item(83, 9).Value = 0
$ws.Cells.Item(83, 10).Value = 3597.6667
$ws.Cells.Item(83, 11).Value = 0
$ws.Cells.Item(83, 12).Value = 32379.0003
$ws.Cells.Item(83, 14).Value = -42363.0003
$ws.Cells.Item(88, 8).Value = 2004.4
$ws.Cells.Item(88, 10).Value = 1882.5
$ws.Cells.Item(88, 12).Value = 1882.5
$ws.Cells.Item(88, 14).Value = -2694.5
$ws.Cells.Item(91, 8).Value = 2004.4
$ws.Cells.Item(91, 10).Value = 1882.5
$ws.Cells.Item(91, 12).Value = 1882.5
$ws.Cells.Item(91, 14).Value = -4690.5
$ws.Cells.Item(100, 8).Value = 3066.2104
$ws.Cells.Item(100, 9).Value = 1104.9166
$ws.Cells.Item(100, 11).Value = 1104.9166
$ws.Cells.Item(100, 13).Value = -563.9166
$ws.Cells.Item(128, 8).Value = 45000
$ws.Cells.Item(128, 10).Value = 45000
$ws.Cells.Item(128, 12).Value = 45000
$ws.Cells.Item(128, 14).Value = -54960
$ws.Cells.Item(132, 8).Value = 2539.8723
$ws.Cells.Item(132, 9).Value = 2287.7908
$ws.Cells.Item(132, 11).Value = 6863.3724
$ws.Cells.Item(132, 13).Value = -4333.3724
$ws.Cells.Item(137, 8).Value = 528382.4
$ws.Cells.Item(137, 10).Value = 2908
$ws.Cells.Item(137, 12).Value = 8724
$ws.Cells.Item(137, 14).Value = -13824
$ws.Cells.Item(138, 8).Value = 4341.552
$ws.Cells.Item(138, 10).Value = 4079.4285
$ws.Cells.Item(138, 12).Value = 12238.2855
$ws.Cells.Item(138, 14).Value = -22518.2855

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 1368.3695
$ws.Cells.Item(32, 9).Value = 467.05264
$ws.Cells.Item(32, 11).Value = 467.05264
$ws.Cells.Item(32, 13).Value = -180.05264
$ws.Cells.Item(61, 8).Value = 8989.799999999999
$ws.Cells.Item(61, 9).Value = 4960.9
$ws.Cells.Item(61, 11).Value = 4960.9
$ws.Cells.Item(61, 13).Value = -4748.9
$ws.Cells.Item(74, 8).Value = 3619.72
$ws.Cells.Item(74, 9).Value = 1228.4117
$ws.Cells.Item(74, 11).Value = 1228.4117
$ws.Cells.Item(74, 13).Value = -354.4117000000001
$ws.Cells.Item(77, 8).Value = 3619.72
$ws.Cells.Item(77, 9).Value = 1228.4117
$ws.Cells.Item(77, 11).Value = 6142.058500000001
$ws.Cells.Item(77, 13).Value = -1774.058500000001
$ws.Cells.Item(102, 8).Value = 1302.75
$ws.Cells.Item(102, 9).Value = 1186.742
$ws.Cells.Item(102, 11).Value = 1186.742
$ws.Cells.Item(102, 13).Value = 435.258
$ws.Cells.Item(132, 8).Value = 1713.8667
$ws.Cells.Item(132, 9).Value = 1361.5853
$ws.Cells.Item(132, 10).Value = 5324.75
$ws.Cells.Item(132, 11).Value = 4084.7559
$ws.Cells.Item(132, 12).Value = 15974.25
$ws.Cells.Item(132, 13).Value = -1554.7559
$ws.Cells.Item(132, 14).Value = -21034.25
$ws.Cells.Item(136, 8).Value = 8989.799999999999
$ws.Cells.Item(136, 9).Value = 4960.9
$ws.Cells.Item(136, 11).Value = 14882.7
$ws.Cells.Item(136, 13).Value = -12332.7

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(52, 8).Value = 0
$ws.Cells.Item(52, 10).Value = 0
$ws.Cells.Item(52, 12).Value = 0
$ws.Cells.Item(80, 8).Value = 799.3
$ws.Cells.Item(80, 10).Value = 681.6667
$ws.Cells.Item(80, 12).Value = 681.6667
$ws.Cells.Item(80, 14).Value = -2677.6667
$ws.Cells.Item(83, 8).Value = 799.3
$ws.Cells.Item(83, 10).Value = 681.6667
$ws.Cells.Item(83, 12).Value = 3408.3335
$ws.Cells.Item(83, 14).Value = -13392.3335
$ws.Cells.Item(92, 8).Value = 149999.5
$ws.Cells.Item(92, 10).Value = 149999.5
$ws.Cells.Item(92, 12).Value = 149999.5
$ws.Cells.Item(92, 14).Value = -154991.5
$ws.Cells.Item(121, 8).Value = 0
$ws.Cells.Item(121, 10).Value = 0
$ws.Cells.Item(121, 12).Value = 0

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(58, 8).Value = 5548.8667
$ws.Cells.Item(58, 9).Value = 3769.5833
$ws.Cells.Item(58, 11).Value = 3769.5833
$ws.Cells.Item(58, 13).Value = -3566.5833
$ws.Cells.Item(74, 8).Value = 35189.332
$ws.Cells.Item(74, 10).Value = 47641.5
$ws.Cells.Item(74, 12).Value = 47641.5
$ws.Cells.Item(74, 14).Value = -49389.5
$ws.Cells.Item(77, 8).Value = 35189.332
$ws.Cells.Item(77, 10).Value = 47641.5
$ws.Cells.Item(77, 12).Value = 142924.5
$ws.Cells.Item(77, 14).Value = -151660.5
$ws.Cells.Item(132, 8).Value = 1889.6842
$ws.Cells.Item(132, 9).Value = 1800.25
$ws.Cells.Item(132, 11).Value = 5400.75
$ws.Cells.Item(132, 13).Value = -2870.75
$ws.Cells.Item(134, 8).Value = 3450.1777
$ws.Cells.Item(134, 9).Value = 2506.6177
$ws.Cells.Item(134, 11).Value = 7519.853099999999
$ws.Cells.Item(134, 13).Value = -4984.853099999999
$ws.Cells.Item(136, 8).Value = 5548.8667
$ws.Cells.Item(136, 9).Value = 3769.5833
$ws.Cells.Item(136, 11).Value = 11308.7499
$ws.Cells.Item(136, 13).Value = -8758.749899999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(4, 8).Value = 663853.4
$ws.Cells.Item(4, 9).Value = 1212507
$ws.Cells.Item(4, 11).Value = 3637521
$ws.Cells.Item(4, 13).Value = -3637409
$ws.Cells.Item(5, 8).Value = 4729.4
$ws.Cells.Item(5, 10).Value = 4479.6
$ws.Cells.Item(5, 12).Value = 13438.8
$ws.Cells.Item(5, 14).Value = -13662.8
$ws.Cells.Item(46, 8).Value = 7839.2
$ws.Cells.Item(46, 9).Value = 300
$ws.Cells.Item(46, 10).Value = 9724
$ws.Cells.Item(46, 11).Value = 900
$ws.Cells.Item(46, 12).Value = 29172
$ws.Cells.Item(46, 13).Value = -809
$ws.Cells.Item(46, 14).Value = -29354
$ws.Cells.Item(58, 8).Value = 7001.5
$ws.Cells.Item(58, 9).Value = 4000
$ws.Cells.Item(58, 11).Value = 12000
$ws.Cells.Item(58, 13).Value = -11872
$ws.Cells.Item(98, 8).Value = 4515.68
$ws.Cells.Item(98, 9).Value = 4722.2856
$ws.Cells.Item(98, 10).Value = 4435.3335
$ws.Cells.Item(98, 11).Value = 14166.8568
$ws.Cells.Item(98, 12).Value = 13306.0005
$ws.Cells.Item(98, 13).Value = -12668.8568
$ws.Cells.Item(98, 14).Value = -16302.0005
$ws.Cells.Item(122, 8).Value = 2509
$ws.Cells.Item(122, 9).Value = 1667.25
$ws.Cells.Item(122, 11).Value = 15005.25
$ws.Cells.Item(122, 13).Value = -12555.25
$ws.Cells.Item(128, 8).Value = 193985
$ws.Cells.Item(128, 9).Value = 193985
$ws.Cells.Item(128, 11).Value = 581955
$ws.Cells.Item(128, 13).Value = -576975
$ws.Cells.Item(131, 8).Value = 3100.8
$ws.Cells.Item(131, 9).Value = 1104.875
$ws.Cells.Item(131, 11).Value = 3314.625
$ws.Cells.Item(131, 13).Value = 1725.375
$ws.Cells.Item(135, 8).Value = 4729.4
$ws.Cells.Item(135, 10).Value = 4479.6
$ws.Cells.Item(135, 12).Value = 40316.4
$ws.Cells.Item(135, 14).Value = -45386.4

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(43, 8).Value = 5933.1113
$ws.Cells.Item(43, 9).Value = 5933.1113
$ws.Cells.Item(43, 10).Value = 0
$ws.Cells.Item(43, 11).Value = 5933.1113
$ws.Cells.Item(43, 12).Value = 0
$ws.Cells.Item(43, 13).Value = -5782.1113
$ws.Cells.Item(132, 8).Value = 4402.2593
$ws.Cells.Item(132, 9).Value = 4598.6
$ws.Cells.Item(132, 10).Value = 1948
$ws.Cells.Item(132, 11).Value = 13795.8
$ws.Cells.Item(132, 12).Value = 5844
$ws.Cells.Item(132, 13).Value = -11265.8
$ws.Cells.Item(132, 14).Value = -10904

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(117, 8).Value = 99000
$ws.Cells.Item(117, 10).Value = 99000
$ws.Cells.Item(117, 12).Value = 99000
$ws.Cells.Item(117, 14).Value = -108178
$ws.Cells.Item(122, 8).Value = 5432.5386
$ws.Cells.Item(122, 9).Value = 2695.5
$ws.Cells.Item(122, 10).Value = 7778.5713
$ws.Cells.Item(122, 11).Value = 8086.5
$ws.Cells.Item(122, 12).Value = 23335.7139
$ws.Cells.Item(122, 13).Value = -5636.5
$ws.Cells.Item(122, 14).Value = -28235.7139
$ws.Cells.Item(136, 8).Value = 1896.0167
$ws.Cells.Item(136, 9).Value = 1190.2195
$ws.Cells.Item(136, 11).Value = 3570.6585
$ws.Cells.Item(136, 13).Value = -1020.6585

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(80, 8).Value = 48599.6
$ws.Cells.Item(80, 10).Value = 48599.6
$ws.Cells.Item(80, 12).Value = 48599.6
$ws.Cells.Item(80, 14).Value = -50595.6
$ws.Cells.Item(83, 8).Value = 48599.6
$ws.Cells.Item(83, 10).Value = 48599.6
$ws.Cells.Item(83, 12).Value = 145798.8
$ws.Cells.Item(83, 14).Value = -155782.8
$ws.Cells.Item(132, 8).Value = 1444.0209
$ws.Cells.Item(132, 9).Value = 1350.3182
$ws.Cells.Item(132, 11).Value = 4050.9546
$ws.Cells.Item(132, 13).Value = -1520.9546
$ws.Cells.Item(136, 8).Value = 9135.1875
$ws.Cells.Item(136, 9).Value = 8300.839
$ws.Cells.Item(136, 11).Value = 24902.517
$ws.Cells.Item(136, 13).Value = -22352.517

# Cells that no longer have a value (removed entirely)
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(80, 13).ClearContents()
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(83, 13).ClearContents()
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(52, 14).ClearContents()
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(121, 14).ClearContents()
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(43, 14).ClearContents()
